$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) for rows being updated so that
# values like "1.00" or "0.995" are not auto-converted to numbers by Excel,
# matching the original inline-string ("text") representation in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.878.74'
$ws.Range("E2").Value = '  +2.20%  '

$ws.Range("D3").Value = '3.394.53'
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").Value = '560.50'
$ws.Range("E5").Value = '  +1.26%  '

$ws.Range("D6").Value = '175.30'
$ws.Range("E6").Value = '  +2.28%  '

$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +2.70%  '

$ws.Range("D8").Value = '3.378.32'
$ws.Range("E8").Value = '  +1.77%  '

$ws.Range("E9").Value = '  +0.25%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +10.45%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '0.631'
$ws.Range("E11").Value = '  +2.32%  '

$ws.Range("D12").Value = '54.62'
$ws.Range("E12").Value = '  +1.63%  '

$ws.Range("D13").Value = '0.0000276'
$ws.Range("E13").Value = '  +4.48%  '

$ws.Range("D14").Value = '9.14'
$ws.Range("E14").Value = '  +2.55%  '

$ws.Range("D15").Value = '3.950.02'
$ws.Range("E15").Value = '  +2.79%  '

$ws.Range("D16").Value = '18.28'
$ws.Range("E16").Value = '  +3.75%  '

$ws.Range("D17").Value = '3.401.72'
$ws.Range("E17").Value = '  +2.48%  '

$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").Value = '11.91'
$ws.Range("E19").Value = '  +2.91%  '

$ws.Range("D20").Value = '64.935.94'
$ws.Range("E20").Value = '  +2.48%  '

$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  +2.58%  '

$ws.Range("D22").Value = '471.96'
$ws.Range("E22").Value = '  +16.77%  '

$ws.Range("D23").Value = '4.98'
$ws.Range("E23").Value = '  +16.72%  '

$ws.Range("D24").Value = '4.14'
$ws.Range("E24").Value = '  +2.64%  '

$ws.Range("D25").Value = '86.62'
$ws.Range("E25").Value = '  +5.21%  '

$ws.Range("D26").Value = '13.69'
$ws.Range("E26").Value = '  +5.63%  '

$ws.Range("D27").Value = '10.87'
$ws.Range("E27").Value = '  +2.07%  '

$ws.Range("D28").Value = '2.87'
$ws.Range("E28").Value = '  +5.20%  '

$ws.Range("D29").Value = '8.84'
$ws.Range("E29").Value = '  +1.89%  '

$ws.Range("D30").Value = '30.68'
$ws.Range("E30").Value = '  +5.65%  '

$ws.Range("D31").Value = '6.72'
$ws.Range("E31").Value = '  +3.52%  '

$ws.Range("D32").Value = '11.53'
$ws.Range("E32").Value = '  +2.19%  '

$ws.Range("D33").Value = '580.02'
$ws.Range("E33").Value = '  -1.36%  '

$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  +2.96%  '

$ws.Range("D35").Value = '60.01'
$ws.Range("E35").Value = '  +4.26%  '

$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").Value = '0.141'
$ws.Range("E37").Value = '  -4.45%  '

$ws.Range("D38").Value = '35.95'
$ws.Range("E38").Value = '  +0.87%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '3.47'
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0756'
$ws.Range("E40").Value = '  +1.96%  '

$ws.Range("D41").Value = '0.373'
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("D42").Value = '3.107.13'
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.36%  '

$ws.Range("D44").Value = '2.87'
$ws.Range("E44").Value = '  +1.94%  '

$ws.Range("D45").Value = '2.52'
$ws.Range("E45").Value = '  +1.87%  '

$ws.Range("D46").Value = '0.0413'
$ws.Range("E46").Value = '  +2.35%  '

$ws.Range("D47").Value = '3.21'
$ws.Range("E47").Value = '  +1.67%  '

$ws.Range("E48").Value = '  +4.79%  '

$ws.Range("D49").Value = '2.57'
$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("D50").Value = '8.40'
$ws.Range("E50").Value = '  +4.90%  '

$ws.Range("D51").Value = '136.63'
$ws.Range("E51").Value = '  +3.71%  '

